$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" column (G) lists the users/systems who recorded each
# session. Some rows show both "dnasr281@gmail.com" and "System" having
# recorded - normalize the ordering so "System" is listed first.
$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
